$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A74").Value = "Neufahrn b. Freising"
$ws.Range("A73").Value = "Mühldorf a. Inn"
$ws.Range("A56").Value = "Weilheim i. OB"
$ws.Range("A45").Value = "Pfaffenhofen a. d. Ilm"
$ws.Range("A44").Value = "Lauf a.d. Pegnitz"
$ws.Range("A31").Value = "Neumarkt i. d. OBp"
$ws.Range("A25").Value = "Weiden"
$ws.Range("A22").Value = "Hof"

$ws.Range("A23").Select()
